# Updates the transition-probability matrix on Sheet1 with recalculated
# values (more games were simulated, changing the underlying counts/ratios).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B2" = 0.2758620689655172
    "C2" = 0.3793103448275862
    "J2" = 0.03448275862068965
    "P2" = 0.1724137931034483
    "S2" = 0.1379310344827586

    "C3" = 0.1666666666666667
    "P3" = 0.8333333333333334

    "J4" = 0.1111111111111111
    "P4" = 0.6666666666666666
    "S4" = 0.2222222222222222

    "J6" = 0.1951219512195122
    "Q6" = 0.1951219512195122
    "R6" = 0.04878048780487805
    "S6" = 0.5609756097560976

    "B7" = 0.04444444444444445
    "F7" = 0.02222222222222222
    "J7" = 0.1111111111111111
    "Q7" = 0.1333333333333333
    "R7" = 0.08888888888888889
    "S7" = 0.6

    "B8" = 0.03478260869565217
    "D8" = 0.03478260869565217
    "F8" = 0.04347826086956522
    "J8" = 0.06956521739130435
    "O8" = 0.01739130434782609
    "Q8" = 0.1043478260869565
    "R8" = 0.05217391304347826
    "S8" = 0.6434782608695652

    "D9" = 0.01694915254237288
    "F9" = 0.1186440677966102
    "J9" = 0.1186440677966102
    "Q9" = 0.1355932203389831
    "R9" = 0.1016949152542373
    "S9" = 0.5084745762711864

    "B10" = 0.05761316872427984
    "D10" = 0.0205761316872428
    "F10" = 0.07407407407407407
    "J10" = 0.1193415637860082
    "O10" = 0.01234567901234568
    "Q10" = 0.2098765432098765
    "R10" = 0.08641975308641975
    "S10" = 0.4197530864197531

    "G11" = 0.1645569620253164
    "J11" = 0.05063291139240506
    "K11" = 0.2278481012658228
    "L11" = 0.4683544303797468
    "S11" = 0.08860759493670886

    "G12" = 0.6756756756756757
    "J12" = 0.1891891891891892
    "S12" = 0.1351351351351351

    "G13" = 0.5833333333333334
    "J13" = 0.25
    "S13" = 0.1666666666666667

    "H15" = 0.2325581395348837
    "I15" = 0.1162790697674419
    "J15" = 0.3488372093023256
    "K15" = 0.06976744186046512
    "M15" = 0.02325581395348837
    "S15" = 0.2093023255813954

    "H16" = 0.2857142857142857
    "I16" = 0.1428571428571428
    "J16" = 0.2857142857142857
    "K16" = 0.04761904761904762
    "S16" = 0.2380952380952381

    "F17" = 0.01204819277108434
    "H17" = 0.1927710843373494
    "I17" = 0.0963855421686747
    "J17" = 0.3734939759036144
    "K17" = 0.06024096385542169
    "M17" = 0.02409638554216868
    "O17" = 0.0963855421686747
    "S17" = 0.144578313253012

    "F18" = 0.05
    "H18" = 0.175
    "I18" = 0.075
    "J18" = 0.4
    "K18" = 0.1
    "O18" = 0.1
    "S18" = 0.1

    "F19" = 0.002898550724637681
    "H19" = 0.2318840579710145
    "I19" = 0.1217391304347826
    "J19" = 0.3043478260869565
    "K19" = 0.1304347826086956
    "M19" = 0.02608695652173913
    "O19" = 0.06086956521739131
    "S19" = 0.1217391304347826
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
